$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextNumberValue($cellRef, $textVal) {
    $ws.Range("H26").NumberFormat = "@"
    $ws.Range("H26").Value = $textVal
    $ws.Range("H26").Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

Set-TextNumberValue "C2" "10300009"
Set-TextNumberValue "C3" "10300015"
Set-TextNumberValue "C4" "10300215"
Set-TextNumberValue "C5" "10300024"
Set-TextNumberValue "C11" "10300129"
Set-TextNumberValue "C12" "10300116"
Set-TextNumberValue "C13" "10300117"
Set-TextNumberValue "C14" "10300121"
Set-TextNumberValue "C24" "137274899"

$ws.Range("H26").ClearContents()
$ws.Range("G2").Copy()
$ws.Range("H26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$v_F2 = @'
Expected condition failed: waiting for visibility of element located by By.id: lblServiceID (tried for 30 second(s) with 500 milliseconds interval)
'@
$ws.Range("F2").Value = $v_F2

$v_F4 = @'
Expected condition failed: waiting for element to be clickable: [[ChromeDriver: chrome on WINDOWS (90906add1415bc7936b4e50a6a63f496)] -> xpath: //a[@id="idTask"]] (tried for 60 second(s) with 500 milliseconds interval)
Build info: version: '3.141.59', revision: 'e82be7d358', time: '2018-11-14T08:17:03'
System info: host: 'SIPL92', ip: '10.212.130.2', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '20'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 122.0.6261.129, chrome: {chromedriverVersion: 122.0.6261.128 (f18a44fedeb..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:49604}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 90906add1415bc7936b4e50a6a63f496
'@
$ws.Range("F4").Value = $v_F4

$v_F11 = @'
no such element: Unable to locate element: {"method":"css selector","selector":"#lblServiceID"}
  (Session info: chrome-headless-shell=122.0.6261.129)
For documentation on this error, please visit: https://www.seleniumhq.org/exceptions/no_such_element.html
Build info: version: '3.141.59', revision: 'e82be7d358', time: '2018-11-14T08:17:03'
System info: host: 'SIPL92', ip: '10.212.130.16', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '20'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome-headless-shell, browserVersion: 122.0.6261.129, chrome: {chromedriverVersion: 122.0.6261.128 (f18a44fedeb..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:63458}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 3e8e4ea869dab0400e39a0f5e6cf0ea0
*** Element info: {Using=id, value=lblServiceID}
'@
$ws.Range("F11").Value = $v_F11

$v_F12 = @'
no such element: Unable to locate element: {"method":"css selector","selector":"#lblServiceID"}
  (Session info: chrome-headless-shell=122.0.6261.129)
For documentation on this error, please visit: https://www.seleniumhq.org/exceptions/no_such_element.html
Build info: version: '3.141.59', revision: 'e82be7d358', time: '2018-11-14T08:17:03'
System info: host: 'SIPL92', ip: '10.212.130.16', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '20'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome-headless-shell, browserVersion: 122.0.6261.129, chrome: {chromedriverVersion: 122.0.6261.128 (f18a44fedeb..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:63458}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 3e8e4ea869dab0400e39a0f5e6cf0ea0
*** Element info: {Using=id, value=lblServiceID}
'@
$ws.Range("F12").Value = $v_F12

$v_F13 = @'
no such element: Unable to locate element: {"method":"css selector","selector":"#lblServiceID"}
  (Session info: chrome-headless-shell=122.0.6261.129)
For documentation on this error, please visit: https://www.seleniumhq.org/exceptions/no_such_element.html
Build info: version: '3.141.59', revision: 'e82be7d358', time: '2018-11-14T08:17:03'
System info: host: 'SIPL92', ip: '10.212.130.16', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '20'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome-headless-shell, browserVersion: 122.0.6261.129, chrome: {chromedriverVersion: 122.0.6261.128 (f18a44fedeb..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:63458}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 3e8e4ea869dab0400e39a0f5e6cf0ea0
*** Element info: {Using=id, value=lblServiceID}
'@
$ws.Range("F13").Value = $v_F13

$v_F14 = @'
no such element: Unable to locate element: {"method":"css selector","selector":"#lblServiceID"}
  (Session info: chrome-headless-shell=122.0.6261.129)
For documentation on this error, please visit: https://www.seleniumhq.org/exceptions/no_such_element.html
Build info: version: '3.141.59', revision: 'e82be7d358', time: '2018-11-14T08:17:03'
System info: host: 'SIPL92', ip: '10.212.130.16', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '20'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome-headless-shell, browserVersion: 122.0.6261.129, chrome: {chromedriverVersion: 122.0.6261.128 (f18a44fedeb..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:63458}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 3e8e4ea869dab0400e39a0f5e6cf0ea0
*** Element info: {Using=id, value=lblServiceID}
'@
$ws.Range("F14").Value = $v_F14

$v_F17 = @'
Parts not available
'@
$ws.Range("F17").Value = $v_F17

$v_E18 = @'
PASS
'@
$ws.Range("E18").Value = $v_E18

$v_F18 = @'
Expected condition failed: waiting for visibility of all elements located by By.xpath: //*[@ng-form="FDXUPSFOrm"] (tried for 60 second(s) with 500 milliseconds interval)
Build info: version: '3.141.59', revision: 'e82be7d358', time: '2018-11-14T08:17:03'
System info: host: 'SIPL92', ip: '10.212.130.2', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '20'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 122.0.6261.129, chrome: {chromedriverVersion: 122.0.6261.128 (f18a44fedeb..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:49604}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 90906add1415bc7936b4e50a6a63f496
'@
$ws.Range("F18").Value = $v_F18

$v_F26 = @'
Cannot invoke "org.openqa.selenium.WebElement.isDisplayed()" because "element" is null
'@
$ws.Range("F26").Value = $v_F26
